# TC09_Bento_MultiFilter_* "startup" sheet holds four Neo4j Cypher queries
# per row 2-4 (column B = "query", column C = "StatQuery", for the
# CasesTab / SamplesTab / FilesTab rows). This edit updates the
# `tp.endocrine_therapy_type` filter in every one of those queries from
# "Other" to "None" to match the newly available data set, per the commit
# message "updated bento scripts as per availability of objects for new
# data set".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldFilter = 'tp.endocrine_therapy_type IN ["Other"]'
$newFilter = 'tp.endocrine_therapy_type IN ["None"]'

foreach ($addr in @("B2", "B3", "B4", "C2", "C3", "C4")) {
    $cell = $ws.Range($addr)
    $cell.Value2 = $cell.Value2.Replace($oldFilter, $newFilter)
}

# The saved workbook shows the selection moved from C4 to B3.
$ws.Range("B3").Select()
